{"js": "// Turn the single long run in each of the 4 paragraphs below into a run made\n// of several <w:t> segments separated by manual line breaks (<w:br/>), at the\n// points described by the diff. We locate each split point with a short,\n// unique snippet of text that straddles the boundary, then rewrite just that\n// snippet so the break (\"\\v\", Word's manual-line-break character) lands\n// exactly between the two halves - this preserves all the surrounding text\n// and the run's formatting.\n\nconst body = context.document.body;\n\nasync function splitAt(straddlingText, breakIndex) {\n  // `straddlingText` is literal text that spans the desired break point;\n  // `breakIndex` is the offset (within that text) where the break goes.\n  const left = straddlingText.slice(0, breakIndex);\n  const right = straddlingText.slice(breakIndex);\n\n  const results = body.search(straddlingText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      \"Expected exactly 1 match for '\" + straddlingText + \"', got \" + results.items.length\n    );\n  }\n\n  results.items[0].insertText(left + \"\\v\" + right, \"Replace\");\n  await context.sync();\n}\n\n// 1) \"Programa\" paragraph: break before \"Realiza\u00e7\u00e3o de atividade experimental\"\nawait splitAt(\n  \"ligas n\u00e3o-ferrosas. Curvas TTT e CCT (TRC). Realiza\u00e7\u00e3o de atividade experimental (8 horas-aula)\",\n  \"ligas n\u00e3o-ferrosas. Curvas TTT e CCT (TRC). \".length\n);\n\n// 2) \"M\u00e9todo\" paragraph: break before \"O aluno ser\u00e1 avaliado...\" and before\n//    \"O desenvolvimento do aluno...\"\nawait splitAt(\n  \"exerc\u00edcios num\u00e9ricos. O aluno ser\u00e1 avaliado ao longo do semestre\",\n  \"exerc\u00edcios num\u00e9ricos. \".length\n);\nawait splitAt(\n  \"atividade experimental (8 horas-aula). O desenvolvimento do aluno ao longo do curso\",\n  \"atividade experimental (8 horas-aula). \".length\n);\n\n// 3) \"Norma de recupera\u00e7\u00e3o\" paragraph: break before \"M\u00e9dia final\"\nawait splitAt(\n  \"0 (zero) a 10 (dez). M\u00e9dia final = (NF + PR) / 2\",\n  \"0 (zero) a 10 (dez). \".length\n);\n\n// 4) \"Bibliografia\" paragraph: break between each of the 10 numbered\n//    references (9 breaks total)\nconst biblioBoundaries = [\n  [\", LTC Editora, 2013.\", \"2. ASKELAND, D.R., P\"],\n  [\"GE, S\u00e3o Paulo, 2008.\", \"3. SHACKELFORD, J.F.\"],\n  [\" ed., Pearson, 2008.\", \"4. GARCIA, A. Solidi\"],\n  [\". Ed. UNICAMP, 2001.\", \"5. READEY, D. W. Kin\"],\n  [\"ress, 1st. Ed. 2016.\", \"6. SHEWMON, P.G. Dif\"],\n  [\" McGraw-Hill, 1963. \", \"7. SHEWMON, P.G. Pha\"],\n  [\" McGraw-Hill, 1969. \", \"8. HUMPHREYS, F.J, H\"],\n  [\"na. Pergamon, 2004. \", \"9. BILLMEYER JR., F.\"],\n  [\"ns, New York, 1984. \", \"10. WILSON, E.A. Wor\"],\n];\n\nfor (const [left, right] of biblioBoundaries) {\n  await splitAt(left + right, left.length);\n}\n", "ps1": "# Turn the single long run in each of the 4 target paragraphs into a run\n# made of several text segments separated by manual line breaks (^l, which\n# Word serializes as <w:br/>), at the points described by the diff.\n#\n# For each split point we Find/Replace a short, unique snippet of text that\n# straddles the desired break (using the whole-document Find so it is not\n# sensitive to which paragraph/run currently holds the text), inserting a\n# manual line break (\"^l\" in ReplaceWith) exactly between the two halves.\n# wdReplaceAll (2) is used, but every snippet is engineered to match exactly\n# once, so this is equivalent to a single targeted replacement.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Once {\n    param([string]$FindText, [string]$ReplaceText)\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    #   MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n    $ok = $find.Execute($FindText, $true, $false, $false, $false, $false, $true, 1, $false, $ReplaceText, 2)\n    if (-not $ok) {\n        throw \"Find/Replace failed for: $FindText\"\n    }\n}\n\n# Programa paragraph: break before \"Realiza\u00e7\u00e3o de atividade experimental\"\nReplace-Once \"ligas n\u00e3o-ferrosas. Curvas TTT e CCT (TRC). Realiza\u00e7\u00e3o de atividade experimental (8 horas-aula)\" \"ligas n\u00e3o-ferrosas. Curvas TTT e CCT (TRC). ^lRealiza\u00e7\u00e3o de atividade experimental (8 horas-aula)\"\n\n# M\u00e9todo paragraph: break before \"O aluno ser\u00e1 avaliado\"\nReplace-Once \"exerc\u00edcios num\u00e9ricos. O aluno ser\u00e1 avaliado ao longo do semestre\" \"exerc\u00edcios num\u00e9ricos. ^lO aluno ser\u00e1 avaliado ao longo do semestre\"\n\n# M\u00e9todo paragraph: break before \"O desenvolvimento do aluno\"\nReplace-Once \"atividade experimental (8 horas-aula). O desenvolvimento do aluno ao longo do curso\" \"atividade experimental (8 horas-aula). ^lO desenvolvimento do aluno ao longo do curso\"\n\n# Norma de recupera\u00e7\u00e3o paragraph: break before \"M\u00e9dia final\"\nReplace-Once \"0 (zero) a 10 (dez). M\u00e9dia final = (NF + PR) / 2\" \"0 (zero) a 10 (dez). ^lM\u00e9dia final = (NF + PR) / 2\"\n\n# Bibliografia paragraph: break before reference #2\nReplace-Once \", LTC Editora, 2013.2. ASKELAND, D.R., P\" \", LTC Editora, 2013.^l2. ASKELAND, D.R., P\"\n\n# Bibliografia paragraph: break before reference #3\nReplace-Once \"GE, S\u00e3o Paulo, 2008.3. SHACKELFORD, J.F.\" \"GE, S\u00e3o Paulo, 2008.^l3. SHACKELFORD, J.F.\"\n\n# Bibliografia paragraph: break before reference #4\nReplace-Once \" ed., Pearson, 2008.4. GARCIA, A. Solidi\" \" ed., Pearson, 2008.^l4. GARCIA, A. Solidi\"\n\n# Bibliografia paragraph: break before reference #5\nReplace-Once \". Ed. UNICAMP, 2001.5. READEY, D. W. Kin\" \". Ed. UNICAMP, 2001.^l5. READEY, D. W. Kin\"\n\n# Bibliografia paragraph: break before reference #6\nReplace-Once \"ress, 1st. Ed. 2016.6. SHEWMON, P.G. Dif\" \"ress, 1st. Ed. 2016.^l6. SHEWMON, P.G. Dif\"\n\n# Bibliografia paragraph: break before reference #7\nReplace-Once \" McGraw-Hill, 1963. 7. SHEWMON, P.G. Pha\" \" McGraw-Hill, 1963. ^l7. SHEWMON, P.G. Pha\"\n\n# Bibliografia paragraph: break before reference #8\nReplace-Once \" McGraw-Hill, 1969. 8. HUMPHREYS, F.J, H\" \" McGraw-Hill, 1969. ^l8. HUMPHREYS, F.J, H\"\n\n# Bibliografia paragraph: break before reference #9\nReplace-Once \"na. Pergamon, 2004. 9. BILLMEYER JR., F.\" \"na. Pergamon, 2004. ^l9. BILLMEYER JR., F.\"\n\n# Bibliografia paragraph: break before reference #10\nReplace-Once \"ns, New York, 1984. 10. WILSON, E.A. Wor\" \"ns, New York, 1984. ^l10. WILSON, E.A. Wor\"\n"}
